$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: A=1, B=2, C=3, D=4, E=5
# Each entry: row, D-value (or $null to skip), E-value (or $null to skip)
$updates = @(
    @(2,  "57.208.67",   "  +4.46%  "),
    @(3,  "2.507.52",    "  +3.24%  "),
    @(4,  $null,         "  -0.01%  "),
    @(5,  "494.82",      "  +3.17%  "),
    @(6,  "153.68",      "  +11.47%  "),
    @(7,  "0.997",       "  +0.06%  "),
    @(8,  $null,         "  +3.22%  "),
    @(9,  "2.523.61",    "  +3.01%  "),
    @(10, $null,         "  +4.62%  "),
    @(11, "5.77",        "  +5.76%  "),
    @(12, $null,         "  +4.58%  "),
    @(13, $null,         "  +1.28%  "),
    @(14, "2.948.32",    "  +3.48%  "),
    @(15, "57.338.30",   $null),
    @(16, "21.38",       "  +4.84%  "),
    @(17, $null,         "  +3.05%  "),
    @(18, "2.528.81",    "  +3.23%  "),
    @(19, $null,         "  +6.25%  "),
    @(20, "10.33",       "  +5.81%  "),
    @(21, "323.51",      "  +3.21%  "),
    @(22, $null,         "  +0.31%  "),
    @(23, "5.93",        "  +5.85%  "),
    @(24, "58.48",       $null),
    @(25, "0.411",       "  +2.20%  "),
    @(26, $null,         "  +0.52%  "),
    @(27, $null,         "  -0.40%  "),
    @(28, "2.619.80",    "  +2.98%  "),
    @(29, "7.63",        "  +4.39%  "),
    @(30, "0.0₃0832",    "  +7.93%  "),
    @(31, $null,         "  +0.17%  "),
    @(32, "151.49",      "  +1.99%  "),
    @(33, $null,         "  +4.51%  "),
    @(34, $null,         "  +2.57%  "),
    @(35, "5.30",        "  +3.31%  "),
    @(36, "1.16",        "  +4.62%  "),
    @(37, "3.82",        "  +6.44%  "),
    @(38, $null,         "  +5.81%  "),
    @(39, $null,         "  +10.02%  "),
    @(40, "34.43",       "  +3.94%  "),
    @(41, "3.55",        $null),
    @(42, "0.622",       "  +4.07%  "),
    @(43, "0.0562",      "  +3.56%  "),
    @(44, "0.994",       "  +0.24%  "),
    @(45, "4.94",        "  +6.68%  "),
    @(46, "270.13",      "  +6.56%  "),
    @(47, $null,         "  +5.01%  "),
    @(48, "0.0230",      "  +4.11%  "),
    @(49, $null,         "  +1.27%  "),
    @(50, "18.03",       "  +6.28%  "),
    @(51, "1.903.04",    "  -1.56%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($null -ne $dVal) {
        $cell = $ws.Cells.Item($row, 4)
        # Force text format so numeric-looking strings (e.g. "494.82",
        # "5.30") are preserved verbatim instead of being coerced into
        # real numbers (which would also drop significant trailing zeros).
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}
